# Update the email addresses on the "AccountCreationData" worksheet and
# move the active selection, matching the author's latest data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreationData")

$ws.Range("A2").Value = "aewtest1@gmail.com"
$ws.Range("A3").Value = "aewtest2@gmail.com"
$ws.Range("A4").Value = "aewtest3@gmail.com"

$ws.Activate()
$ws.Range("E17").Select()
